# Apply updated crypto price/volume figures (inline-string text cells).
# Numeric-looking "D" values must be forced to Text format first, otherwise
# the COM layer auto-converts them to floating point numbers (losing the
# original text formatting, e.g. trailing zeros / exact decimal digits).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.856.06"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").Value = "1.891.07"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7757"
$ws.Range("E5").Value = "  -3.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.84"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3145"
$ws.Range("E8").Value = "  -2.94%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07414"
$ws.Range("E9").Value = "  +2.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.33"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08135"
$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("E12").Value = "  -2.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.483"
$ws.Range("E13").Value = "  +1.72%  "

$ws.Range("D14").Value = "1.887.56"
$ws.Range("E14").Value = "  -1.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.28"
$ws.Range("E15").Value = "  -1.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.186"
$ws.Range("E16").Value = "  +2.33%  "

$ws.Range("D17").Value = "29.923.91"
$ws.Range("E17").Value = "  -0.95%  "

$ws.Range("E18").Value = "  -1.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.78"
$ws.Range("E19").Value = "  -1.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007860"
$ws.Range("E20").Value = "  +0.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.107"
$ws.Range("E22").Value = "  -1.18%  "

$ws.Range("D23").Value = "2.133.04"
$ws.Range("E23").Value = "  -2.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1587"
$ws.Range("E25").Value = "  -2.67%  "

$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.54"
$ws.Range("E27").Value = "  -2.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.80"
$ws.Range("E28").Value = "  -0.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.040"
$ws.Range("E29").Value = "  -4.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.450"
$ws.Range("E30").Value = "  +4.38%  "

$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.504"
$ws.Range("E32").Value = "  -0.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.097"
$ws.Range("E33").Value = "  -1.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05570"
$ws.Range("E34").Value = "  -2.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.247"
$ws.Range("E35").Value = "  -3.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7614"
$ws.Range("E36").Value = "  +1.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.009"
$ws.Range("E37").Value = "  +0.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.646"
$ws.Range("E38").Value = "  -3.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01934"
$ws.Range("E39").Value = "  -1.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.789"
$ws.Range("E40").Value = "  -0.81%  "

$ws.Range("D41").Value = "1.165.86"
$ws.Range("E41").Value = "  +12.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4468"
$ws.Range("E42").Value = "  -0.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "74.04"
$ws.Range("E43").Value = "  +0.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.024"
$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8523"
$ws.Range("E45").Value = "  -0.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.900"
$ws.Range("E47").Value = "  -1.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.32"
$ws.Range("E48").Value = "  -0.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.907"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.086"
$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.526"
$ws.Range("E51").Value = "  -1.09%  "
